$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row to reflect the new, simplified beneficiary import layout.
$ws.Range("B1").Value = "Account Number*"
$ws.Range("C1").Value = "Beneficiary Name (Last, First)"
$ws.Range("D1").Value = "Beneficiary Type (Primary / Contingent)"
$ws.Range("E1").Value = "Beneficiary Allocation (%)"
$ws.Range("F1").Value = "Beneficiary Relationship (Spouse, Child, etc.)"
$ws.Range("A1").Value = "Client ID"

# Remove the now-obsolete Contingent-specific columns.
$ws.Range("G1:H1").Clear()

# Update selection / view to match the new layout.
$ws.Range("B8").Select()
